$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "RF Predictetion" -> "RF Prediction"
$ws.Range("C1").Value = "RF Prediction"

# Replace row 16 data: socfb-A-anon -> socfb-B-anon, with updated prediction results
$ws.Range("A16").Value = "socfb-B-anon"
$ws.Range("B16").Value = "Facebook Networks"
$ws.Range("C16").Value = "Social Networks"
$ws.Range("D16").Value = "Social Networks"
$ws.Range("E16").Value = 0.71809776051547003
$ws.Range("F16").Value = "Facebook Networks"
$ws.Range("G16").Value = 0.37855867805294202

# Widen column A and update the saved selection
$ws.Columns("A").ColumnWidth = 23.67
$ws.Range("J33").Select()
